# "fix create cv excel file"
#
# - Drop the unused Sheet2 / Sheet3 tabs.
# - Remove all review comments from Sheet1 (B1, C1, T1, U1) together with
#   their backing VML/author data.
# - Rename the "Annual Salary" column header (Q1) to "Monthly Salary".
# - Update the saved view state (scroll position / active selection).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove Sheet2 and Sheet3 ---------------------------------------------
foreach ($name in @("Sheet2", "Sheet3")) {
    try {
        $wb.Worksheets.Item($name).Delete()
    } catch {
        # sheet already absent - nothing to do
    }
}

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# --- Remove every cell comment on Sheet1 -----------------------------------
while ($ws1.Comments.Count -gt 0) {
    $ws1.Comments.Item(1).Delete()
}

# --- Rename the salary column header ---------------------------------------
$ws1.Range("Q1").Value = "Monthly Salary"

# --- Update the view: scroll to column O, select T6 ------------------------
$ws1.Range("T6").Select()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1

Write-Output "Edit applied"
